# Daily attendance processing - 2025-11-10 18:54:34
# Rotate the "Recorded By" (column G) comma-separated list of names/emails
# so that the last entry moves to the front (right rotation), for every
# row that has more than one recorder listed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $text = [string]$val
    if ($text -notmatch ",") { continue }

    $parts = $text -split ",\s*"
    if ($parts.Count -lt 2) { continue }

    $rotated = @($parts[-1]) + $parts[0..($parts.Count - 2)]
    $cell.Value2 = [string]::Join(", ", $rotated)
}
